$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:F4").NumberFormat = "@"

$ws.Range("E2").Value = "2019-01-01"
$ws.Range("F2").Value = "2019-01-03"
$ws.Range("E3").Value = "2019-01-02"
$ws.Range("F3").Value = "2019-01-04"
$ws.Range("E4").Value = "2019-01-03"
$ws.Range("F4").Value = "2019-01-05"

$ws.Range("C2:C4").NumberFormat = "0"

$ws.Range("C3:C4").Select() | Out-Null

$ws.Range("C1:C4").EntireColumn.AutoFit() | Out-Null
